# full_system.pptx edit:
#  - widen/enlarge the "Inst ROM" box (Rounded Rectangle 59) and bump its
#    label font size 12 -> 14
#  - remove the two duplicate "Inst ROM" boxes (Rounded Rectangle 64 / 41)
#  - nudge the two core-label textboxes ("Multi-Core" / "SIMD-Core") down
#  - add two new "Inst ROM" boxes (clones of the resized one) under the
#    Multi-Core and SIMD-Core columns

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Resize "Inst ROM" (Rounded Rectangle 59) and grow its text to 14pt
$instRom = $s.Shapes.Item("Rounded Rectangle 59")
$instRom.Left = 54
$instRom.Top = 480
$instRom.Width = 78
$instRom.Height = 30
$instRom.TextFrame.TextRange.Font.Size = 14

# 2) Delete the two redundant "Inst ROM" boxes, leaving their connectors
#    (now unglued on that end) in place
$s.Shapes.Item("Rounded Rectangle 64").Delete()
$s.Shapes.Item("Rounded Rectangle 41").Delete()

# 3) Move the "Multi-Core" / "SIMD-Core" textboxes down a bit
$s.Shapes.Item("TextBox 85").Top = 36
$s.Shapes.Item("TextBox 86").Top = 36

# 4) Add two new "Inst ROM" boxes (same look as the resized one) under the
#    Multi-Core and SIMD-Core stacks
$instRom55 = $instRom.Duplicate()
$instRom55.Name = "Rounded Rectangle 55"
$instRom55.Left = 198
$instRom55.Top = 480
$instRom55.Width = 78
$instRom55.Height = 30

$instRom57 = $instRom.Duplicate()
$instRom57.Name = "Rounded Rectangle 57"
$instRom57.Left = 486
$instRom57.Top = 480
$instRom57.Width = 78
$instRom57.Height = 30
